$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures (D = Price, E = Volume(1h)) for rows 2-51.
# 'Numeric' rows are those whose new Price text looks like a plain number (e.g. '1.00',
# '0.999') -- Excel would silently coerce those to a numeric value and drop the exact
# text formatting, so those cells get forced to Text format before the value is set.
$updates = @{
    2 = @{ D='62.658.72'; DNumeric=$false; E='  -4.12%  ' }
    3 = @{ D='3.049.04'; DNumeric=$false; E='  -3.32%  ' }
    4 = @{ E='  +0.16%  ' }
    5 = @{ D='543.51'; DNumeric=$true; E='  -4.65%  ' }
    6 = @{ D='133.74'; DNumeric=$true; E='  -11.70%  ' }
    7 = @{ D='1.00'; DNumeric=$true; E='  +0.06%  ' }
    8 = @{ D='3.044.44'; DNumeric=$false; E='  -3.22%  ' }
    9 = @{ D='0.488'; DNumeric=$true; E='  -3.04%  ' }
    10 = @{ D='0.154'; DNumeric=$true; E='  -4.43%  ' }
    11 = @{ D='6.35'; DNumeric=$true; E='  -11.43%  ' }
    12 = @{ D='0.457'; DNumeric=$true; E='  -3.15%  ' }
    13 = @{ D='34.61'; DNumeric=$true; E='  -5.66%  ' }
    14 = @{ D='0.0000213'; DNumeric=$true; E='  -6.14%  ' }
    15 = @{ D='3.545.88'; DNumeric=$false; E='  -3.04%  ' }
    16 = @{ D='62.774.67'; DNumeric=$false; E='  -3.92%  ' }
    17 = @{ E='  -2.51%  ' }
    18 = @{ D='3.052.68'; DNumeric=$false; E='  -3.09%  ' }
    19 = @{ D='6.60'; DNumeric=$true; E='  -3.85%  ' }
    20 = @{ D='479.55'; DNumeric=$true; E='  -10.95%  ' }
    21 = @{ D='13.31'; DNumeric=$true; E='  -5.37%  ' }
    22 = @{ D='0.696'; DNumeric=$true; E='  -2.94%  ' }
    23 = @{ D='7.02'; DNumeric=$true; E='  -6.98%  ' }
    24 = @{ D='76.98'; DNumeric=$true; E='  -3.32%  ' }
    25 = @{ D='12.11'; DNumeric=$true; E='  -7.60%  ' }
    26 = @{ E='  -0.16%  ' }
    27 = @{ D='8.25'; DNumeric=$true; E='  -9.34%  ' }
    28 = @{ D='2.69'; DNumeric=$true; E='  -5.19%  ' }
    29 = @{ D='1.00'; DNumeric=$true; E='  +0.00%  ' }
    30 = @{ D='1.92'; DNumeric=$true; E='  -11.91%  ' }
    31 = @{ D='26.00'; DNumeric=$true; E='  -2.18%  ' }
    32 = @{ D='1.12'; DNumeric=$true; E='  -3.79%  ' }
    33 = @{ E='  -8.41%  ' }
    34 = @{ D='57.82'; DNumeric=$true; E='  +8.63%  ' }
    35 = @{ D='508.38'; DNumeric=$true; E='  -8.74%  ' }
    36 = @{ D='5.93'; DNumeric=$true; E='  -4.19%  ' }
    37 = @{ D='5.05'; DNumeric=$true; E='  -8.83%  ' }
    38 = @{ D='0.0396'; DNumeric=$true; E='  -12.86%  ' }
    39 = @{ D='3.075.64'; DNumeric=$false; E='  +0.20%  ' }
    40 = @{ D='0.0784'; DNumeric=$true; E='  -6.20%  ' }
    41 = @{ E='  -4.64%  ' }
    42 = @{ D='8.01'; DNumeric=$true; E='  -4.64%  ' }
    43 = @{ D='2.57'; DNumeric=$true; E='  -13.14%  ' }
    44 = @{ D='0.250'; DNumeric=$true; E='  -5.09%  ' }
    45 = @{ E='  +0.04%  ' }
    46 = @{ D='2.02'; DNumeric=$true; E='  -10.51%  ' }
    47 = @{ D='120.19'; DNumeric=$true; E='  -0.10%  ' }
    48 = @{ D='24.30'; DNumeric=$true; E='  -4.97%  ' }
    49 = @{ D='0.106'; DNumeric=$true; E='  -4.01%  ' }
    50 = @{ D='2.37'; DNumeric=$true; E='  +58.30%  ' }
    51 = @{ D='0.0₃0489'; DNumeric=$false; E='  -8.75%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey('D')) {
        $priceCell = $ws.Cells.Item($row, 4)
        if ($vals['DNumeric']) {
            $priceCell.NumberFormat = '@'
        }
        $priceCell.Value = $vals['D']
    }
    if ($vals.ContainsKey('E')) {
        $ws.Cells.Item($row, 5).Value = $vals['E']
    }
}